$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new play-by-play yardage entries ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "12 4 3 -2 3 10 11 14 1 8 7 12 3 2 11 3 1 4 3 -1 0 0 -1 3 8 5 7 -1 0 2 3 10 0 0 8 2 0 8 1 -3 1 8 6 0 3 4 1 0 5 2 11 1 3 11 5 10 3 3 3 2 0 2 1 3 -1 5 2 2 8 2 3 2 4 -3 10 2 5 -3 1 2 -3 3 3 13 5 2 0 13 1 2 -3 4 17 10 1 6 1 3 9 0 3 1 4 0 8 4 3 17 0 5 2 1 4 1 1 3 -2 3 6 3 1 7 3 13 1 2 7 2 2 2 2 1 3 11 2 2 3 5 3 2 2 4 1 -1 9 0 5 2 18 2 8 -2 3 3 14 5 0 5 3 1 1 0 4 1 5 0 2 2 1 3 4 4 5 1 2 1 0 1 3 -4 -1 2 6 14 3 2 3 5 -3 4 1 -5 1 6 2 4 7 17 1 3 -3 3 0 12 2 1 0 2 5 5 0 -1 0 7 5 -2 4 19 15 1 11 7 1 27 4 12 4 2 1 2 7 0 6 1 2 6 2 5 2 7 39 17 6 5 1 5 11 0 4 2 0 4 0 3 1 0 2 3 0 14 0 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3 4 0 0 7 2 7 1 3 2 3 3 6 0 0 2 -1 5 4 8 2 1 3"
$ws.Range("C2").Value = "13 7 4 6 6 3 4 11 4 9 23 1 4 14 18 7 3 3 1 6 2 8 -6 -1 5 -3 9 -1 6 6 8 2 7 3 2 10 3 8 10 4 2 8 4 6 0 3 8 7 2 5 12 -3 15 2 -1 3 4 0 4 6 -3 2 4 2 6 11 4 4 1 1 6 -1 7 4 6 5 3 6 2 8 6 6 -2 20 2 5 6 8 10 2 1 5 1 3 10 -1 1 -1 5 3 -3 2 7 3 5 4 9 3 7 8 0 6 1 1 1 2 4 3 3 6 5 2 8 3 7 5 2 2 3 8 6 0 16 3 6 2 3 9 5 3 9 3 10 -5 4 -2 2 4 9 3 9 2 2 1 -2 6 2 5 2 7 2 13 6 6 3 3 6 20 9 2 2 6 4 6 6 6 3 6 6 -2 1 11 8 2 6 -1 1 5 3 2 1 9 2 9 11 9 9 14 2 4 1 3 8 -1 3 1 6 -1 3 2 44 2 4 2 2 3 3 5 4 -3 4 2 -3 3 2 1 9 9 2 2 2 1 0 3 1 3 15 22 12 7 2 3 1 1 1 21 3 2 1 2 14 17 2 1 5 13 -1 7 -1 5 12 2 10 0 4 1 6 2 7 0 6 0 5 6 3 9 7 3 1 8 8 1 4 4 4 3 2 12 7 9 7 3 -1 10 2 10 4 4 3 4 2 2 6 6 5 4 3 1 3 4 9 -4 1 9 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4 2 5 5 3 2 2 17 2 3 6 0 3 10 0 5 6 12 8 13 13 1 3 5 1 1 16 3 3 2 6 2 2"
$ws.Range("B3").Value = "16 1 11 12 7 3 7 9 9 5 18 5 12 6 9 7 6 8 7 5 0 9 5 9 4 5 23 24 4 8 24 7 6 1 14 8 3 8 14 4 23 12 8 7 22 5 7 8 1 2 4 10 2 6 3 8 19 3 4 9 20 4 6 12 14 3 4 1 12 3 -1 10 26 4 6 14 4 1 1 28 25 4 13 15 7 15 10 5 12 42 8 8 17 6 12 15 5 7 14 14 6 7 8 7 5 21 0 3 16 -2 14 18 22 16 2 13 5 22 17 5 17 6 10 6 12 0 13 17 6 7 5 12 9 2 39 4 15 5 7 -3 23 11 7 2 13 12 6 13 15 5 39 9 6 4 4 26 49 4 11 12 6 11 23 28 4 12 3 7 7 15 11 2 2 18 1 -3 2 9 17 5 8 15 10 34 7 39 13 6 0 8 7 5 16 7 20 7 6 3 12 49 16 13 -4 9 64 17 18 17 14 0 8 13 21 9 6 12 1 9 8 9 9 7 0 10 0 6 5 11 4 19 16 17 4 6 9 7 18 13 7 0 4 7 15 10 20 19 5 12 8 14 9 7 0 6 6 12 19 29 7 10 9 8 9 4 -3 14 2 19 6 9 4 16 6 20 6 17 8 8 8 2 7 30 -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23 7 7 14 5 0 13 20 10 21 49 2 4 49 5 14 3 4 5 4"
$ws.Range("C3").Value = "10 11 2 18 13 28 8 4 19 3 9 -3 1 9 6 10 12 9 -1 5 8 25 7 19 9 0 23 20 25 20 16 16 9 5 8 1 9 5 26 17 21 3 8 1 2 9 10 26 5 12 6 14 8 38 19 8 18 12 9 7 4 16 0 13 13 11 10 12 3 7 5 16 4 9 10 4 1 6 9 3 19 -1 21 33 16 20 11 2 -2 11 23 10 17 24 19 4 30 11 14 1 13 3 6 14 4 5 -6 7 23 27 8 15 11 6 10 20 -1 9 17 5 2 10 10 0 7 2 15 9 7 19 5 8 11 5 14 6 6 9 18 27 9 11 1 7 10 17 16 4 4 17 -2 13 10 23 7 9 21 4 11 1 11 -1 9 2 19 5 14 7 19 11 20 4 15 -1 14 19 9 7 6 5 9 25 5 7 22 8 8 37 11 13 5 9 12 6 6 18 5 7 4 9 8 16 9 14 9 19 7 22 11 6 5 5 17 6 9 9 19 12 6 9 5 18 6 10 10 6 6 15 3 19 2 10 1 11 5 9 8 11 24 -1 21 -2 8 9 14 5 14 -5 21 7 11 26 20 3 8 15 5 11 0 6 3 15 7 4 3 25 36 3 0 15 6 17 12 27 1 9 1 7 11 5 7 18 1 5 11 14 6 6 6 36 5 3 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11 6 25 6 28 3 36 9 9 30 3 15 13 21 15 3 4 8 1"

# --- OFF sheet: weekly offense totals (Week 15 logged + Week 16 simulated) ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 185
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 32
$ws.Range("G2").Value = 33
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 18
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 19
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 118
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 118
$ws.Range("G3").Value = 31
$ws.Range("I3").Value = 71
$ws.Range("J3").Value = 47
$ws.Range("L3").Value = 263
$ws.Range("M3").Value = 175
$ws.Range("Q3").Value = 469

# --- DEF sheet: weekly defense totals ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 169
$ws.Range("D2").Value = 9
$ws.Range("F2").Value = 48
$ws.Range("G2").Value = 67
$ws.Range("J2").Value = 32
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 160
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 91
$ws.Range("G3").Value = 26
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 47
$ws.Range("L3").Value = 285
$ws.Range("M3").Value = 192
$ws.Range("Q3").Value = 509

# --- ST sheet: special teams totals + per-game arrays ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 54
$ws.Range("D2").Value = 56
$ws.Range("F2").Value = 86
$ws.Range("G2").Value = 82
$ws.Range("J2").Value = 57
$ws.Range("K2").Value = 55
$ws.Range("L2").Value = 27
$ws.Range("M2").Value = 20
$ws.Range("B3").Value = 29
$ws.Range("D3").Value = "35 52 62 31 50 49 30 54 40 33 57 42 39 37 58 51 46 51 25 45 24 32 50 38 38 55 51 38 46 30 41 52 50 51 41 60 66 51 50 52 54 62 32 55 54 32 51 43 40 35 53 44 46 64 53 46"
$ws.Range("B4").Value = "66 67 57 66 52 60 65 66 64 59 60 64 53 69 45 52 67 66 64 63 63 66 64 66 64"
$ws.Range("D4").Value = "0 2 11 0 7 -1 0 0 0 0 0 0 0 0 19 17 0 6 0 0 0 0 0 16 0 0 0 0 0 0 0 17 7 26 3 16 0 10 0 0 17 0 0 0 18 0 11 1 0 0 6 0 0 0 10 -2"
$ws.Range("B5").Value = "18 22 16 24 19 28 23 1 27 65 13 24 18 27 16 21 29 29 26 19 26 28 0 25 15"
$ws.Range("D5").Value = "0 0 0 0 0 0 23 0"
$ws.Range("B6").Value = "25 4 13 19 14 45 27 24 10 16 24 32 22 29 16 24 21 25 19 21 25 28 21 16 17 27 18 24 22 26 19 30 38 20 26 33 24"

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("E3").Value = 3

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 24
$ws.Range("D4").Value = 10
